# Automatic update of files.
# - Update B29 (Taxonsorteringsordning) to 56575
# - Append a new observation row (row 30) with the new record's data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 29: bump the sort-order value
$ws.Range("B29").Value = 56575

# New row 30 data
$ws.Range("A30").Value = 111982668
$ws.Range("B30").Value = 90800
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 4364
$ws.Range("F30").Value = "Dropptaggsvamp"
$ws.Range("G30").Value = "Hydnellum ferrugineum"
$ws.Range("H30").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I30").Value = "'3"
$ws.Range("J30").Value = "fruktkroppar"
$ws.Range("K30").Value = "'"

$ws.Range("P30").Value = "Oxögat, Boh"
$ws.Range("Q30").Value = 306298
$ws.Range("R30").Value = 6525575
$ws.Range("S30").Value = 50
$ws.Range("T30").Value = "Västra Götaland"
$ws.Range("U30").Value = "Tanum"
$ws.Range("V30").Value = "Bohuslän"
$ws.Range("W30").Value = "Naverstad"

$ws.Range("Y30").Value = "'2023-09-09"
$ws.Range("Z30").Value = "13:36"
$ws.Range("AA30").Value = "'2023-09-09"
$ws.Range("AB30").Value = "13:36"

$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AT30").Value = "'"

$ws.Range("AW30").Value = "Mattias Drejby"
$ws.Range("AX30").Value = "Mattias Drejby"
$ws.Range("AY30").Value = "'"
